# Updated interval in gss
# Re-run of the golden-section-search / BFGS trace: the (X,Y,F,Lambda,DX,DY,DF,
# Gradient,Hk) log on Sheet1 now starts from a different initial interval, so
# rows 2-9 get new numbers and the search now takes until row 16 to converge
# (7 new rows, 10-16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column A formatting (style index 1, same as A2:A9) down to the new rows
$ws.Range("A9").Copy()
$ws.Range("A10:A16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = -0.6914893617021276
$ws.Range("E2").Value = 6.332875812655695
$ws.Range("F2").Value = 2.162561938137853
$ws.Range("G2").Value = 1.224384658965657
$ws.Range("H2").Value = -2.339012326833017
$ws.Range("I2").Value = '[-0.34148182 -0.19333786]'
$ws.Range("J2").Value = '[[1 0]
 [0 1]]'

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2.162561938137853
$ws.Range("C3").Value = 1.224384658965657
$ws.Range("D3").Value = -3.030501688535145
$ws.Range("E3").Value = 0.3695944430204291
$ws.Range("F3").Value = -0.374416873583963
$ws.Range("G3").Value = 0.349913111403884
$ws.Range("H3").Value = -0.3797285322287296
$ws.Range("I3").Value = '[ 0.70846967 -1.25132707]'
$ws.Range("J3").Value = '[[ 0.43893511 -0.56106489]
 [-0.56106489  0.43893511]]'

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1.78814506455389
$ws.Range("C4").Value = 1.574297770369541
$ws.Range("D4").Value = -3.410230220763874
$ws.Range("E4").Value = [double]"4.348389486982408e-06"
$ws.Range("F4").Value = [double]"-7.539340009632411e-06"
$ws.Range("G4").Value = [double]"-7.210590500461578e-06"
$ws.Range("H4").Value = [double]"1.64710352410502e-05"
$ws.Range("I4").Value = '[-1.07971966 -1.15532225]'
$ws.Range("J4").Value = '[[-0.25883241 -1.25883241]
 [-1.25883241 -0.25883241]]'

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 1.78813752521388
$ws.Range("C5").Value = 1.57429055977904
$ws.Range("D5").Value = -3.410213749728633
$ws.Range("E5").Value = 0.01670991455705744
$ws.Range("F5").Value = 0.3885634745955935
$ws.Range("G5").Value = 0.3898262869501039
$ws.Range("H5").Value = -0.4611932377656336
$ws.Range("I5").Value = '[-1.07974867 -1.15532131]'
$ws.Range("J5").Value = '[[10.92081691  9.92081691]
 [ 9.92081691 10.92081691]]'

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 2.176700999809474
$ws.Range("C6").Value = 1.964116846729144
$ws.Range("D6").Value = -3.871406987494267
$ws.Range("E6").Value = 0.338101536010981
$ws.Range("F6").Value = -0.3328567624763523
$ws.Range("G6").Value = 0.3521682144404832
$ws.Range("H6").Value = -0.3477684534712639
$ws.Range("I6").Value = '[ 1.01456732 -1.01152516]'
$ws.Range("J6").Value = '[[-8.88762967 -9.88762967]
 [-9.88762967 -8.88762967]]'

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 1.843844237333121
$ws.Range("C7").Value = 2.316285061169627
$ws.Range("D7").Value = -4.219175440965531
$ws.Range("E7").Value = 0.01325720590178662
$ws.Range("F7").Value = 0.2509446763804477
$ws.Range("G7").Value = 0.2503026655927907
$ws.Range("H7").Value = -0.2173433848841384
$ws.Range("I7").Value = '[-0.88291935 -0.83449205]'
$ws.Range("J7").Value = '[[11.50767869 10.50767869]
 [10.50767869 11.50767869]]'

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 2.094788913713569
$ws.Range("C8").Value = 2.566587726762418
$ws.Range("D8").Value = -4.436518825849669
$ws.Range("E8").Value = 0.2529903630436824
$ws.Range("F8").Value = -0.1319169042482222
$ws.Range("G8").Value = 0.125997523658159
$ws.Range("H8").Value = -0.06601048217374572
$ws.Range("I8").Value = '[ 0.50916697 -0.51029648]'
$ws.Range("J8").Value = '[[ -9.85748566 -10.85748566]
 [-10.85748566  -9.85748566]]'

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 1.962872009465347
$ws.Range("C9").Value = 2.692585250420577
$ws.Range("D9").Value = -4.502529308023415
$ws.Range("E9").Value = 0.005894366301292123
$ws.Range("F9").Value = 0.04764076506913217
$ws.Range("G9").Value = 0.04769363866475862
$ws.Range("H9").Value = -0.009322861658427328
$ws.Range("I9").Value = '[-0.19094954 -0.19991973]'
$ws.Range("J9").Value = '[[21.18954848 20.18954848]
 [20.18954848 21.18954848]]'

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 2.010512774534479
$ws.Range("C10").Value = 2.740278889085336
$ws.Range("D10").Value = -4.511852169681842
$ws.Range("E10").Value = 0.2385519770820237
$ws.Range("F10").Value = -0.01289629486605715
$ws.Range("G10").Value = 0.01234762559289049
$ws.Range("H10").Value = -0.0006676234920650614
$ws.Range("I10").Value = '[ 0.05288253 -0.05293894]'
$ws.Range("J10").Value = '[[-19.88613026 -20.88613026]
 [-20.88613026 -19.88613026]]'

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 1.997616479668422
$ws.Range("C11").Value = 2.752626514678226
$ws.Range("D11").Value = -4.512519793173907
$ws.Range("E11").Value = 0.002148989501692404
$ws.Range("F11").Value = 0.002924244183692304
$ws.Range("G11").Value = 0.002925379642346293
$ws.Range("H11").Value = [double]"-3.564106019915414e-05"
$ws.Range("I11").Value = '[-0.01191378 -0.01244215]'
$ws.Range("J11").Value = '[[56.38032323 55.38032323]
 [55.38032323 56.38032323]]'

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 2.000540723852114
$ws.Range("C12").Value = 2.755551894320572
$ws.Range("D12").Value = -4.512555434234106
$ws.Range("E12").Value = 0.1947057166132569
$ws.Range("F12").Value = -0.0006970853783498931
$ws.Range("G12").Value = 0.0003568218426592829
$ws.Range("H12").Value = [double]"-1.424713099495989e-06"
$ws.Range("I12").Value = '[ 0.00269855 -0.00271427]'
$ws.Range("J12").Value = '[[-55.10018693 -56.10018693]
 [-56.10018693 -55.10018693]]'

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 1.999843638473764
$ws.Range("C13").Value = 2.755908716163232
$ws.Range("D13").Value = -4.512556858947206
$ws.Range("E13").Value = 0.00200162891401638
$ws.Range("F13").Value = 0.0002764971372095637
$ws.Range("G13").Value = 0.0002779861710324738
$ws.Range("H13").Value = [double]"-3.197547808753143e-07"
$ws.Range("I13").Value = '[-0.00078019 -0.0015241 ]'
$ws.Range("J13").Value = '[[60.60853959 59.60853959]
 [59.60853959 60.60853959]]'

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 2.000120135610974
$ws.Range("C14").Value = 2.756186702334264
$ws.Range("D14").Value = -4.512557178701987
$ws.Range("E14").Value = 0.2480648326997222
$ws.Range("F14").Value = -0.0001066291382532114
$ws.Range("G14").Value = 0.0001900286047393784
$ws.Range("H14").Value = [double]"-8.86314239778585e-08"
$ws.Range("I14").Value = '[ 0.00059935 -0.00059654]'
$ws.Range("J14").Value = '[[-59.33023254 -60.33023254]
 [-60.33023254 -59.33023254]]'

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 2.000013506472721
$ws.Range("C15").Value = 2.756376730939003
$ws.Range("D15").Value = -4.512557267333411
$ws.Range("E15").Value = [double]"4.348389486982408e-06"
$ws.Range("F15").Value = [double]"6.90268567637986e-07"
$ws.Range("G15").Value = [double]"6.903969236304874e-07"
$ws.Range("H15").Value = [double]"7.462919171530302e-11"
$ws.Range("I15").Value = '[6.73760548e-05 3.78579879e-05]'
$ws.Range("J15").Value = '[[-1508.09888684 -1509.09888684]
 [-1509.09888684 -1508.09888684]]'

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 2.000014196741288
$ws.Range("C16").Value = 2.756377421335927
$ws.Range("D16").Value = -4.512557267258781
$ws.Range("E16").Value = [double]"7.80286833249941e-05"
$ws.Range("F16").Value = [double]"-1.306781742904661e-05"
$ws.Range("G16").Value = [double]"-1.306542536738675e-05"
$ws.Range("H16").Value = [double]"-7.392859657784356e-10"
$ws.Range("I16").Value = '[7.08193832e-05 4.01631984e-05]'
$ws.Range("J16").Value = '[[1509.37827784 1508.37827784]
 [1508.37827784 1509.37827784]]'
